# Generate Report for Handback
# Refreshes the "Latest HO Xliff Generate Date" / handoff / handback
# timestamps that get stamped whenever the handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 ("Latest HO Xliff Generate Date") and de-de!H2 ("Correspond
# Handoff Datetime") both held the same timestamp string -> bump together.
$wsOverview.Range("G2").Value = "2016-08-28 09:04:25"
$wsDeDe.Range("H2").Value = "2016-08-28 09:04:25"

# zh-cn!H2 ("Correspond Handoff Datetime") / zh-cn!K2 ("Correspond Handback
# DateTime")
$wsZhCn.Range("H2").Value = "2016-08-28 09:04:21"
$wsZhCn.Range("K2").Value = "2016-08-28 09:04:37"

# de-de!K2 ("Correspond Handback DateTime")
$wsDeDe.Range("K2").Value = "2016-08-28 09:04:44"
